$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 updates
$ws.Range("D16").Value = "image_20250807110238_ppp0.jpg"

$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "641,530,687,575"
$ws.Range("I16").Style = "Normal"

# Row 17 updates
$ws.Range("D17").Value = "image_20250807110238_ppp0.jpg"

$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "793,481,831,527"
$ws.Range("I17").Style = "Normal"

$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "0.71"
$ws.Range("J17").Style = "Normal"
